$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A16").Value = "Kun je 4 dozen schroeven bestellen?"
$ws.Range("B16").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$ws.Range("C16").Value = "Hoi Johan, `nZou je 4 dozen schroeven kunnen bestellen voor de werkplaats?`nDank je wel!`nGroet, `nRick`nSent using {0}"
$ws.Range("D16").Value = "Bestelling / Levering"
$ws.Range("E16").Value = "Bedankt voor je bericht. Ik neem dit z.s.m. in behandeling."
$ws.Range("F16").Value = "2025-06-26 21:11:46"
$ws.Range("G16").Value = "Ja"
$ws.Range("H16").Value = "Ja"
$ws.Range("I16").Value = "Nee"
$ws.Rows.Item(16).AutoFit()

# Extend conditional formatting ranges to include the new row 16
$dFcs = $ws.Range("D2:D15").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D16"))
}

$gFcs = $ws.Range("G2:G15").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($ws.Range("G2:G16"))
}

$hFcs = $ws.Range("H2:H15").FormatConditions
for ($i = 1; $i -le $hFcs.Count; $i++) {
    $hFcs.Item($i).ModifyAppliesToRange($ws.Range("H2:H16"))
}

$iFcs = $ws.Range("I2:I15").FormatConditions
for ($i = 1; $i -le $iFcs.Count; $i++) {
    $iFcs.Item($i).ModifyAppliesToRange($ws.Range("I2:I16"))
}

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 11
